$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.925.15'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.665.03'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.57'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E6").Value = '  +4.62%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.25'
$ws.Range("E10").Value = '  +3.11%  '
$ws.Range("E11").Value = '  +3.81%  '
$ws.Range("D12").Value = '1.900.09'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '1.658.71'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.12'
$ws.Range("E16").Value = '  +1.81%  '
$ws.Range("D17").Value = '26.912.07'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '234.18'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.01'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.35'
$ws.Range("E22").Value = '  -1.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.21'
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.20'
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.12'
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.88'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0496'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("E32").Value = '  +2.15%  '
$ws.Range("D33").Value = '1.456.89'
$ws.Range("E33").Value = '  -4.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("E34").Value = '  +2.78%  '
$ws.Range("E35").Value = '  +2.89%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.581'
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0168'
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.71'
$ws.Range("E40").Value = '  -3.96%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("E43").Value = '  +6.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.83'
$ws.Range("E44").Value = '  -0.81%  '
$ws.Range("D45").Value = '1.809.61'
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.781'
$ws.Range("E46").Value = '  +0.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.37'
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = [string]::Concat('0.0', [char]0x2086, '0105')
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("E50").Value = '  +3.93%  '
$ws.Range("E51").Value = '  +0.20%  '
